$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 274 (shifts existing rows 274-332 down to 275-333)
$ws.Rows.Item(274).Insert()

# Populate the newly inserted row with the Agra, India colo data
$ws.Range("A274").Value2 = "AGR"
$ws.Range("B274").Value2 = "Agra, India"
$ws.Range("C274").Value2 = "Asia Pacific"
$ws.Range("D274").Value2 = "Agra"
$ws.Range("E274").Value2 = "India"
$ws.Range("F274").Value2 = "IN"
$ws.Range("G274").Value2 = 27.202556
$ws.Range("H274").Value2 = 77.85740699999999

# Match the formatting used by the rest of column A (bold/centered/bordered)
$ws.Range("A275").Copy()
$ws.Range("A274").PasteSpecial(-4122)
$excel.CutCopyMode = 0
